$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Kampagnendaten 2022 für das Sternbild Bootes Konstellation",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Kampagnendaten 2022 für das Bootes Konstellation",
    2
)
